# "Generate Report for Archive"
#
# The localization-status report was regenerated. The three rows that track
# the in-flight files (1f75d576…, cdb13a05…, d647d643…) came back out in a
# different order, and two of them (cdb13a05…, d647d643…) picked up a new
# "In Translation" status on the Overview sheet while their per-locale
# status stays "Ready for handoff".
#
# New row order (rows 3-5) on every sheet: cdb13a05…, d647d643…, 1f75d576…

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay($ws, $addr, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet: columns A (File Name), B (Path And Name), C (Extension),
# D (Publish URL), E (zh-cn status), F (de-de status), G (Latest HO Xliff
# Generate Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "cdb13a05-bab5-4edb-a4bb-685164cf6771.md"
$wsOverview.Range("B3").Value = "e2e\cdb13a05-bab5-4edb-a4bb-685164cf6771.md"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("G3").Value = "2016-08-30 12:50:16"

$wsOverview.Range("A4").Value = "d647d643-a35b-4847-b6d8-24a33984b8bf.md"
$wsOverview.Range("B4").Value = "e2e\d647d643-a35b-4847-b6d8-24a33984b8bf.md"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"
$wsOverview.Range("G4").Value = "2016-08-30 12:50:16"

$wsOverview.Range("A5").Value = "1f75d576-7a84-4e75-91b5-41ad7cde3d93.md"
$wsOverview.Range("B5").Value = "e2e\1f75d576-7a84-4e75-91b5-41ad7cde3d93.md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-30 12:48:42"

Set-HyperlinkDisplay $wsOverview '$B$3' "e2e\cdb13a05-bab5-4edb-a4bb-685164cf6771.md"
Set-HyperlinkDisplay $wsOverview '$B$4' "e2e\d647d643-a35b-4847-b6d8-24a33984b8bf.md"
Set-HyperlinkDisplay $wsOverview '$B$5' "e2e\1f75d576-7a84-4e75-91b5-41ad7cde3d93.md"

# ---------------------------------------------------------------------
# zh-cn sheet: only columns A (Source File Name), G (Latest Handoff File)
# and H (Latest Handoff Datetime) differ between the three rows; every
# other column is identical across them, so only those need updating.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = "cdb13a05-bab5-4edb-a4bb-685164cf6771.md"
$wsZh.Range("G3").Value = "cdb13a05-bab5-4edb-a4bb-685164cf6771.1728faed9e65750bf6bf7d1691f621f0a476881c.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-30 12:49:59"

$wsZh.Range("A4").Value = "d647d643-a35b-4847-b6d8-24a33984b8bf.md"
$wsZh.Range("G4").Value = "d647d643-a35b-4847-b6d8-24a33984b8bf.2686b2ee432072f8429ba7c17cd8e8a0ec71fa60.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-08-30 12:49:59"

$wsZh.Range("A5").Value = "1f75d576-7a84-4e75-91b5-41ad7cde3d93.md"
$wsZh.Range("G5").Value = "1f75d576-7a84-4e75-91b5-41ad7cde3d93.6f168f45fe5e3961929053a2d377cb2813fa9eed.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-08-30 12:48:37"

Set-HyperlinkDisplay $wsZh '$A$3' "cdb13a05-bab5-4edb-a4bb-685164cf6771.md"
Set-HyperlinkDisplay $wsZh '$A$4' "d647d643-a35b-4847-b6d8-24a33984b8bf.md"
Set-HyperlinkDisplay $wsZh '$A$5' "1f75d576-7a84-4e75-91b5-41ad7cde3d93.md"

# ---------------------------------------------------------------------
# de-de sheet: same shape as zh-cn.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = "cdb13a05-bab5-4edb-a4bb-685164cf6771.md"
$wsDe.Range("G3").Value = "cdb13a05-bab5-4edb-a4bb-685164cf6771.1728faed9e65750bf6bf7d1691f621f0a476881c.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-30 12:50:16"

$wsDe.Range("A4").Value = "d647d643-a35b-4847-b6d8-24a33984b8bf.md"
$wsDe.Range("G4").Value = "d647d643-a35b-4847-b6d8-24a33984b8bf.2686b2ee432072f8429ba7c17cd8e8a0ec71fa60.de-de.xlf"
$wsDe.Range("H4").Value = "2016-08-30 12:50:16"

$wsDe.Range("A5").Value = "1f75d576-7a84-4e75-91b5-41ad7cde3d93.md"
$wsDe.Range("G5").Value = "1f75d576-7a84-4e75-91b5-41ad7cde3d93.6f168f45fe5e3961929053a2d377cb2813fa9eed.de-de.xlf"
$wsDe.Range("H5").Value = "2016-08-30 12:48:42"

Set-HyperlinkDisplay $wsDe '$A$3' "cdb13a05-bab5-4edb-a4bb-685164cf6771.md"
Set-HyperlinkDisplay $wsDe '$A$4' "d647d643-a35b-4847-b6d8-24a33984b8bf.md"
Set-HyperlinkDisplay $wsDe '$A$5' "1f75d576-7a84-4e75-91b5-41ad7cde3d93.md"
